$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where column D already carries the "empty placeholder" style (s="3"):
# just fill D with the "Rien pour le moment…" text and add a new E cell
# ("ignore") copying D's formatting so the new cell lands on style 3 too.
$rowsSimple = @(11, 13, 45, 46, 59)
foreach ($r in $rowsSimple) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)

    $dCell.Copy()
    $eCell.PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $dCell.Value = "Rien pour le moment…"
    $eCell.Value = "ignore"
}

# Row 12: column D currently uses the "vertical-center" style (s="2") instead
# of the plain style (s="3") used by its siblings, so pull formatting from a
# neighbour (D11) before writing the value, then handle E12 as above.
$d12 = $ws.Cells.Item(12, 4)
$e12 = $ws.Cells.Item(12, 5)

$ws.Cells.Item(11, 4).Copy()
$d12.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$d12.Value = "Rien pour le moment…"

$d12.Copy()
$e12.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$e12.Value = "ignore"

# Selection / view mirrors the saved workbook: scrolled back to the top,
# E11:E13 selected with E11 active.
$ws.Range("E11:E13").Select() | Out-Null
$ws.Application.ActiveWindow.ScrollRow = 1
